$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this pushes the existing row 4
# ("Buiten Vlaanderen en Brussel" / gebiedscode 99999) down to row 5,
# and (matching Excel's default insert behaviour) copies the formatting
# of the row above into the freshly inserted row.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the "Niet te lokaliseren" entry.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 99993
$ws.Range("C4").Value = "Niet te lokaliseren"
$ws.Range("D4").Value = "Niet te lokaliseren"

# The old row 4 data is now in row 5; bump its volgnr from 3 to 4.
$ws.Range("A5").Value = 4
